$wb = $excel.ActiveWorkbook

# --- Update selection on the "all" sheet before adding the new sheet ---
$allSheet = $wb.Worksheets.Item("all")
[void]$allSheet.Range("B29").Select()

# --- Add the new "robust" sheet after the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "robust"
[void]$ws.Activate()

# --- Column widths ---
$ws.Columns.Item(4).ColumnWidth = 1.625

# --- Row labels first (this establishes shared-string creation order) ---
$labels = "cif", "qcif", "brate", "gray", "fps", "5fps", "1 degree", "2 degree", "3 degree"
$bvals = 22, 18, 8, 9, 30, 54, 18, 22, 24
$cvals = 97, 81, 93, 83, 98, 80, 85, 86, $null

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}

# --- Column headers (SUJ / KAIST) ---
$ws.Range("B4").Value = "SUJ"
$ws.Range("C4").Value = "KAIST"
$ws.Range("E4").Value = "SUJ"
$ws.Range("F4").Value = "KAIST"

# --- Title rows (added last) ---
$ws.Range("A1").Value = "False Negatives"
$ws.Range("A2").Value = "Total comparisons"
$ws.Range("C2").Value = 48400

# --- Remaining numeric data + formulas ---
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 2).Value = $bvals[$i]
    if ($cvals[$i] -ne $null) {
        $ws.Cells.Item($row, 3).Value = $cvals[$i]
    }
    $ws.Cells.Item($row, 5).Formula = "=(B$row-2)/(C2)"
    $ws.Cells.Item($row, 6).Formula = "=(C$row-2)/(C2)"
}

[void]$ws.Range("H9").Select()
